$d = $word.ActiveDocument

$phrases = @("المرحلة الأولى", "المرحلة الثانية", "المرحلة الثالثة")

foreach ($phrase in $phrases) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Text = $phrase
    $find.Replacement.Text = $phrase
    $find.Replacement.Highlight = 4
    $find.Execute([ref]$phrase, $false, $false, $false, $false, $false, $true, 1, $false, $phrase, 2) | Out-Null
}
